$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_6_6_24"
$ws.Range("B2").Value = 0.261187389609702
$ws.Range("C2").Value = -0.7094970306718658
$ws.Range("D2").Value = 0.0974110662837453
$ws.Range("E2").Value = -0.03778898289354804
$ws.Range("F2").Value = 0.8176478743553162
$ws.Range("G2").Value = 1.121598839759827
$ws.Range("H2").Value = 1.334409952163696
$ws.Range("I2").Value = 1.221745491027832

$ws.Range("A3").Value = "model_6_6_23"
$ws.Range("B3").Value = 0.2626522905978483
$ws.Range("C3").Value = -0.6967094132133735
$ws.Range("D3").Value = 0.09904440297645645
$ws.Range("E3").Value = -0.03304994010551021
$ws.Range("F3").Value = 0.8160266280174255
$ws.Range("G3").Value = 1.113209009170532
$ws.Range("H3").Value = 1.331995368003845
$ws.Range("I3").Value = 1.216166377067566

$ws.Range("A4").Value = "model_6_6_19"
$ws.Range("B4").Value = 0.2638078520893516
$ws.Range("C4").Value = -0.6860438494014527
$ws.Range("D4").Value = 0.1009420516084276
$ws.Range("E4").Value = -0.02878243881068054
$ws.Range("F4").Value = 0.8147478103637695
$ws.Range("G4").Value = 1.106211304664612
$ws.Range("H4").Value = 1.329189658164978
$ws.Range("I4").Value = 1.211142539978027

$ws.Range("A5").Value = "model_6_6_22"
$ws.Range("B5").Value = 0.2646009146409307
$ws.Range("C5").Value = -0.6862405635753832
$ws.Range("D5").Value = 0.1041922205818138
$ws.Range("E5").Value = -0.02691942198848229
$ws.Range("F5").Value = 0.8138700723648071
$ws.Range("G5").Value = 1.106340408325195
$ws.Range("H5").Value = 1.324384570121765
$ws.Range("I5").Value = 1.208949208259583

$ws.Range("A6").Value = "model_6_6_20"
$ws.Range("B6").Value = 0.2649180742668088
$ws.Range("C6").Value = -0.676846230287782
$ws.Range("D6").Value = 0.1018838306880262
$ws.Range("E6").Value = -0.02551172158786419
$ws.Range("F6").Value = 0.8135190606117249
$ws.Range("G6").Value = 1.100176692008972
$ws.Range("H6").Value = 1.327797293663025
$ws.Range("I6").Value = 1.207291960716248

$ws.Range("A7").Value = "model_6_6_21"
$ws.Range("B7").Value = 0.2652408007370126
$ws.Range("C7").Value = -0.6831741420904085
$ws.Range("D7").Value = 0.1063303563379089
$ws.Range("E7").Value = -0.02475073820771345
$ws.Range("F7").Value = 0.8131619095802307
$ws.Range("G7").Value = 1.104328513145447
$ws.Range("H7").Value = 1.321223378181458
$ws.Range("I7").Value = 1.206396102905273

$ws.Range("A8").Value = "model_6_6_18"
$ws.Range("B8").Value = 0.2673964460267421
$ws.Range("C8").Value = -0.6636608667224024
$ws.Range("D8").Value = 0.1090051601006946
$ws.Range("E8").Value = -0.01741279534692119
$ws.Range("F8").Value = 0.8107762932777405
$ws.Range("G8").Value = 1.091525793075562
$ws.Range("H8").Value = 1.317269086837769
$ws.Range("I8").Value = 1.197757482528687

$ws.Range("A9").Value = "model_6_6_17"
$ws.Range("B9").Value = 0.2681706057599671
$ws.Range("C9").Value = -0.6542011992552668
$ws.Range("D9").Value = 0.1088352113275788
$ws.Range("E9").Value = -0.01472228879635762
$ws.Range("F9").Value = 0.8099194765090942
$ws.Range("G9").Value = 1.08531928062439
$ws.Range("H9").Value = 1.317520260810852
$ws.Range("I9").Value = 1.194590210914612

$ws.Range("A10").Value = "model_6_6_16"
$ws.Range("B10").Value = 0.2701866768171776
$ws.Range("C10").Value = -0.6393893754842033
$ws.Range("D10").Value = 0.1131615433383344
$ws.Range("E10").Value = -0.007795094598643804
$ws.Range("F10").Value = 0.807688295841217
$ws.Range("G10").Value = 1.07560133934021
$ws.Range("H10").Value = 1.311124205589294
$ws.Range("I10").Value = 1.186434984207153

$ws.Range("A11").Value = "model_6_6_15"
$ws.Range("B11").Value = 0.270249474191689
$ws.Range("C11").Value = -0.6311160369138764
$ws.Range("D11").Value = 0.1096988850306212
$ws.Range("E11").Value = -0.00740109182666937
$ws.Range("F11").Value = 0.8076187968254089
$ws.Range("G11").Value = 1.070173263549805
$ws.Range("H11").Value = 1.316243290901184
$ws.Range("I11").Value = 1.185971140861511

$ws.Range("A12").Value = "model_6_6_14"
$ws.Range("B12").Value = 0.2703742039371297
$ws.Range("C12").Value = -0.629493838949569
$ws.Range("D12").Value = 0.1096225851477152
$ws.Range("E12").Value = -0.006967054408021589
$ws.Range("F12").Value = 0.8074807524681091
$ws.Range("G12").Value = 1.069108843803406
$ws.Range("H12").Value = 1.316356301307678
$ws.Range("I12").Value = 1.185460090637207

$ws.Range("A13").Value = "model_6_6_12"
$ws.Range("B13").Value = 0.2721988795698285
$ws.Range("C13").Value = -0.6125257273303286
$ws.Range("D13").Value = 0.1122901185090787
$ws.Range("E13").Value = -0.0003847824012024326
$ws.Range("F13").Value = 0.8054613471031189
$ws.Range("G13").Value = 1.057976126670837
$ws.Range("H13").Value = 1.31241238117218
$ws.Range("I13").Value = 1.177711129188538

$ws.Range("A14").Value = "model_6_6_13"
$ws.Range("B14").Value = 0.2735029071811721
$ws.Range("C14").Value = -0.606081241081059
$ws.Range("D14").Value = 0.1150592779574428
$ws.Range("E14").Value = 0.003153067095462592
$ws.Range("F14").Value = 0.8040181994438171
$ws.Range("G14").Value = 1.053747773170471
$ws.Range("H14").Value = 1.308318376541138
$ws.Range("I14").Value = 1.173546195030212

$ws.Range("A15").Value = "model_6_6_10"
$ws.Range("B15").Value = 0.2751792675918482
$ws.Range("C15").Value = -0.5825086413462224
$ws.Range("D15").Value = 0.1150201049895384
$ws.Range("E15").Value = 0.01008570230242889
$ws.Range("F15").Value = 0.8021630048751831
$ws.Range("G15").Value = 1.038281917572021
$ws.Range("H15").Value = 1.308376550674438
$ws.Range("I15").Value = 1.165384650230408

$ws.Range("A16").Value = "model_6_6_11"
$ws.Range("B16").Value = 0.2753724065695222
$ws.Range("C16").Value = -0.5848257133548802
$ws.Range("D16").Value = 0.1163575590774074
$ws.Range("E16").Value = 0.01019152102792065
$ws.Range("F16").Value = 0.8019492030143738
$ws.Range("G16").Value = 1.039802074432373
$ws.Range("H16").Value = 1.306398987770081
$ws.Range("I16").Value = 1.165259957313538

$ws.Range("A17").Value = "model_6_6_9"
$ws.Range("B17").Value = 0.2793419589180398
$ws.Range("C17").Value = -0.5519992147507926
$ws.Range("D17").Value = 0.1219608684224608
$ws.Range("E17").Value = 0.02318843446298524
$ws.Range("F17").Value = 0.7975561022758484
$ws.Range("G17").Value = 1.018264770507812
$ws.Range("H17").Value = 1.298115015029907
$ws.Range("I17").Value = 1.149959325790405

$ws.Range("A18").Value = "model_6_6_8"
$ws.Range("B18").Value = 0.2804704855431043
$ws.Range("C18").Value = -0.5390074639415512
$ws.Range("D18").Value = 0.1220544880203492
$ws.Range("E18").Value = 0.02707771471873166
$ws.Range("F18").Value = 0.7963071465492249
$ws.Range("G18").Value = 1.009740829467773
$ws.Range("H18").Value = 1.297976493835449
$ws.Range("I18").Value = 1.145380735397339

$ws.Range("A19").Value = "model_6_6_7"
$ws.Range("B19").Value = 0.2822417081194114
$ws.Range("C19").Value = -0.5176216640504399
$ws.Range("D19").Value = 0.1217026230056468
$ws.Range("E19").Value = 0.03317911786602212
$ws.Range("F19").Value = 0.7943469285964966
$ws.Range("G19").Value = 0.9957096576690674
$ws.Range("H19").Value = 1.298496723175049
$ws.Range("I19").Value = 1.138197660446167

$ws.Range("A20").Value = "model_6_6_6"
$ws.Range("B20").Value = 0.2859212863213911
$ws.Range("C20").Value = -0.4791073317678594
$ws.Range("D20").Value = 0.1238467217767782
$ws.Range("E20").Value = 0.04580973344247619
$ws.Range("F20").Value = 0.7902747392654419
$ws.Range("G20").Value = 0.9704404473304749
$ws.Range("H20").Value = 1.295326828956604
$ws.Range("I20").Value = 1.12332820892334

$ws.Range("A21").Value = "model_6_6_5"
$ws.Range("B21").Value = 0.2906256779699956
$ws.Range("C21").Value = -0.4137773084865946
$ws.Range("D21").Value = 0.1218281599539779
$ws.Range("E21").Value = 0.06389223783514575
$ws.Range("F21").Value = 0.7850683927536011
$ws.Range("G21").Value = 0.9275775551795959
$ws.Range("H21").Value = 1.298311233520508
$ws.Range("I21").Value = 1.102040410041809

$ws.Range("A22").Value = "model_6_6_4"
$ws.Range("B22").Value = 0.294727649509721
$ws.Range("C22").Value = -0.3677195609796813
$ws.Range("D22").Value = 0.1209835595680633
$ws.Range("E22").Value = 0.07698244715386693
$ws.Range("F22").Value = 0.7805286049842834
$ws.Range("G22").Value = 0.8973591327667236
$ws.Range("H22").Value = 1.299559831619263
$ws.Range("I22").Value = 1.086629986763

$ws.Range("A23").Value = "model_6_6_3"
$ws.Range("B23").Value = 0.3013369239656286
$ws.Range("C23").Value = -0.2787927524565632
$ws.Range("D23").Value = 0.1156011034600163
$ws.Range("E23").Value = 0.1000391062622251
$ws.Range("F23").Value = 0.7732141613960266
$ws.Range("G23").Value = 0.8390142917633057
$ws.Range("H23").Value = 1.307517409324646
$ws.Range("I23").Value = 1.059486269950867

$ws.Range("A24").Value = "model_6_6_2"
$ws.Range("B24").Value = 0.3049633473269423
$ws.Range("C24").Value = -0.227177544377271
$ws.Range("D24").Value = 0.1112899117997808
$ws.Range("E24").Value = 0.1127205306939467
$ws.Range("F24").Value = 0.7692008018493652
$ws.Range("G24").Value = 0.8051497340202332
$ws.Range("H24").Value = 1.313891291618347
$ws.Range("I24").Value = 1.044556856155396

$ws.Range("A25").Value = "model_6_6_1"
$ws.Range("B25").Value = 0.312829728598811
$ws.Range("C25").Value = 0.05861218002310553
$ws.Range("D25").Value = 0.1133925095918258
$ws.Range("E25").Value = 0.1982838863063635
$ws.Range("F25").Value = 0.7604950070381165
$ws.Range("G25").Value = 0.6176433563232422
$ws.Range("H25").Value = 1.310782670974731
$ws.Range("I25").Value = 0.9438267946243286

$ws.Range("A26").Value = "model_6_6_0"
$ws.Range("B26").Value = 0.3378535972445872
$ws.Range("C26").Value = 0.1868508324083674
$ws.Range("D26").Value = 0.1979665168909164
$ws.Range("E26").Value = 0.2861014829158793
$ws.Range("F26").Value = 0.7328009605407715
$ws.Range("G26").Value = 0.5335061550140381
$ws.Range("H26").Value = 1.185746192932129
$ws.Range("I26").Value = 0.8404428362846375
